# pas-730: range of fixes for daily run
# The VIN upload template had two extra columns ("STAT" and "CHOICE_TIER")
# that are no longer part of the expected layout. Remove them (as whole
# columns, so everything to their right shifts left), matching the manual
# fix that was made directly in Excel (select column Z, ctrl-click column
# AC, then delete both columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the right-most extra column first (CHOICE_TIER, column AC) and
# then the left-most one (STAT, column Z) so the column letters used below
# still refer to the right cells at the time each delete happens.
$ws.Range("AC1").EntireColumn.Delete()
$ws.Range("Z1").EntireColumn.Delete()

# Reproduce the resulting on-screen selection as closely as possible: both
# deleted columns (now collapsed to Z:Z and AC:AC) ended up selected, with
# the active cell sitting in the second (AC) area.
$r1 = $ws.Range("Z1:Z1048576")
$r2 = $ws.Range("AC1:AC1048576")
$unionRange = $excel.Union($r1, $r2)
$unionRange.Select()
$ws.Range("AC1").Select()

$win = $excel.ActiveWindow
$win.ScrollColumn = 23
$win.ScrollRow = 1
